$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 8: this shifts the existing rows 8..29
# down to 9..30, matching every row-to-row-below data move visible in the
# diff (row 8 -> 9, row 9 -> 10, ... row 29 -> 30).
$ws.Rows.Item(8).Insert()

# Populate the brand-new row 8 with its own record (a new daily price
# observation that didn't exist before).
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44811
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112043
$ws.Range("G8").Value = "Pepino dulce"
$ws.Range("H8").Value = "Cultivar IV Región"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("N8").Value = "$/bandeja 18 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 806
$ws.Range("Q8").Value = 18
$ws.Range("R8").Value = "Hortaliza"
